# OpenTBS demo workbook: enhanced Excel example.
# Adds a "Score" column (with a bordered input cell and a SUM total) to the
# Example #1 block, inserts a new bullet about picture placeholders, and
# rewords/extends the "do not use a formula" bullet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Insert a new row (16) for a bullet that was pushed down one slot ---
# This shifts the old rows 17-30 down to 18-31, turning the old blank
# separator row (18) into row 19, the header row (19) into row 20, and the
# example data row (20) into row 21.
$ws.Rows("16:16").Insert()

# Row 13 (unchanged position): reword the "do not use a formula" bullet.
$ws.Range("B13").Value = '* Do not use a formula in a cell that may have its position changed after the merge (for example under a TBS block). Otherwise Excel will raise an error message.'

# Row 14 (unchanged position): new explanatory continuation line.
$ws.Range("B14").Value = '    This is because the location of formulas are saved a second time in another sub-file for the order of evaluation.'

# Row 15 (unchanged position): now holds the "reference not rearranged"
# bullet that used to sit on row 14.
$ws.Range("B15").Value = '* If a formula uses a reference to a cell that has moved during the merge, then the reference will not be arraged to be the new cell reference. '

# Row 16 (brand new row): the old "cannot change picture" bullet, now last
# in the list. Formatting was already inherited from row 15 by the insert
# (same bullet style), so only the text needs to be set.
$ws.Range("B16").Value = '* You cannot change picture using "ope=changepic". This is because drawing information are not saved directly in the sheet.'

# --- 2. Fill in the "Total:" label + SUM formula on the (already blank) ---
# --- separator row (19) above the Example #1 data table.               ---
$ws.Range("D19").Value = "Total:"
$ws.Range("D19").HorizontalAlignment = -4152

$ws.Range("B20").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Formula = "=SUM(E21:E2000)"
$ws.Range("E19").Font.Bold = $true
$ws.Range("E19").NumberFormat = "#,##0.0"
$excel.CutCopyMode = $false

# --- 3. Add the "Score" column to the Example #1 header/data rows. ---
$ws.Range("D20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "Score"
$excel.CutCopyMode = $false

$ws.Range("D21").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").Value = "[a.score;ope=xlsxNum]"
$ws.Range("E21").NumberFormat = "#,##0.0"
$ws.Range("E21").HorizontalAlignment = -4152
$excel.CutCopyMode = $false

# --- 4. Selection matches the authored file (cursor left on the new ---
# --- "Score" header cell).                                          ---
$ws.Range("E20").Select()
